$d = $word.ActiveDocument
$x = [char]0xD7

function Replace-Text($old, $new) {
    return $d.Content.Find.Execute($old, $true, $false, $false, $false, $false, $true, 1, $false, $new, 2)
}

$pairs = @(
    @("2023-09-17 Sunday", "2023-09-18 Monday"),
    @(("33{0}81=" -f $x), ("83{0}31=" -f $x)),
    @(("84{0}51=" -f $x), ("17{0}59=" -f $x)),
    @(("23{0}96=" -f $x), ("27{0}52=" -f $x)),
    @(("45{0}72=" -f $x), ("47{0}37=" -f $x)),
    @(("49{0}42=" -f $x), ("80{0}78=" -f $x)),
    @(("86{0}79=" -f $x), ("83{0}17=" -f $x)),
    @(("62{0}35=" -f $x), ("56{0}81=" -f $x)),
    @(("99{0}23=" -f $x), ("59{0}21=" -f $x)),
    @(("32{0}96=" -f $x), ("25{0}25=" -f $x)),
    @(("64{0}28=" -f $x), ("56{0}50=" -f $x)),
    @(("80{0}31=" -f $x), ("69{0}97=" -f $x)),
    @(("12{0}48=" -f $x), ("80{0}95=" -f $x)),
    @(("84{0}34=" -f $x), ("65{0}54=" -f $x)),
    @(("92{0}51=" -f $x), ("80{0}80=" -f $x)),
    @(("83{0}24=" -f $x), ("50{0}80=" -f $x)),
    @(("37{0}99=" -f $x), ("19{0}82=" -f $x)),
    @(("14{0}55=" -f $x), ("51{0}23=" -f $x)),
    @(("38{0}58=" -f $x), ("65{0}18=" -f $x)),
    @(("51{0}36=" -f $x), ("20{0}97=" -f $x)),
    @(("22{0}14=" -f $x), ("32{0}38=" -f $x)),
    @(("79{0}19=" -f $x), ("44{0}80=" -f $x)),
    @(("34{0}48=" -f $x), ("73{0}45=" -f $x)),
    @(("43{0}98=" -f $x), ("72{0}66=" -f $x)),
    @(("30{0}37=" -f $x), ("86{0}99=" -f $x)),
    @(("65{0}43=" -f $x), ("33{0}70=" -f $x))
)

$failures = 0
foreach ($pair in $pairs) {
    $old = $pair[0]
    $new = $pair[1]
    $result = Replace-Text $old $new
    if (-not $result) {
        $failures = $failures + 1
        Write-Host "FAILED: $old -> $new"
    }
}
Write-Host "Done. Failures: $failures"
